$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 8
$ws.Cells.Item(2,3).Value = 2
$ws.Cells.Item(2,4).Value = 6
$ws.Cells.Item(2,5).Value = 6
$ws.Cells.Item(2,6).Value = -2
$ws.Cells.Item(2,7).Value = 4
$ws.Cells.Item(2,8).Value = 32
$ws.Cells.Item(2,9).Value = 5
$ws.Cells.Item(2,10).Value = "train_dim2_1"

# Row 3
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 7
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = 6
$ws.Cells.Item(3,5).Value = 7
$ws.Cells.Item(3,6).Value = -1
$ws.Cells.Item(3,7).Value = 5
$ws.Cells.Item(3,8).Value = 21
$ws.Cells.Item(3,9).Value = 5
$ws.Cells.Item(3,10).Value = "train_dim2_1"

# Row 4
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 6
$ws.Cells.Item(4,3).Value = 3
$ws.Cells.Item(4,4).Value = 1
$ws.Cells.Item(4,5).Value = 4
$ws.Cells.Item(4,6).Value = -5
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 65
$ws.Cells.Item(4,9).Value = 5
$ws.Cells.Item(4,10).Value = "train_dim2_1"

# Row 5
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 9
$ws.Cells.Item(5,3).Value = 4
$ws.Cells.Item(5,4).Value = 5
$ws.Cells.Item(5,5).Value = 6
$ws.Cells.Item(5,6).Value = -4
$ws.Cells.Item(5,7).Value = 2
$ws.Cells.Item(5,8).Value = 54
$ws.Cells.Item(5,9).Value = 5
$ws.Cells.Item(5,10).Value = "train_dim2_1"

# Row 6 (new row)
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 7
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 4
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = -3
$ws.Cells.Item(6,7).Value = 3
$ws.Cells.Item(6,8).Value = 43
$ws.Cells.Item(6,9).Value = 5
$ws.Cells.Item(6,10).Value = "train_dim2_1"

# Update selection to match the new active cell shown in the diff
$ws.Range("I1").Select()
